# Applies the commit's row reshuffle + new-match append to Sheet1 of the
# Serie A 2023-2024 odds workbook.
#
# Context: most "changes" in the diff are really just two (or three)
# existing match rows trading places with each other -- columns A:E
# (index / pais / torneio / temporada / data_partida) stay put, while
# F:V (the match-specific data: teams, scores, odds, timestamps, url)
# swap between the rows. On top of that, one brand-new match (row 182,
# Bologna vs Genoa) is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param([int]$RowA, [int]$RowB)

    $rngA = $ws.Range("F$RowA`:V$RowA")
    $rngB = $ws.Range("F$RowB`:V$RowB")

    $valA = $rngA.Value()
    $valB = $rngB.Value()

    $rngA.Value = $valB
    $rngB.Value = $valA
}

# Simple two-row swaps (F:V only; A:E -- same match date -- stay as-is).
Swap-MatchRows 8 9
Swap-MatchRows 24 25
Swap-MatchRows 36 37
Swap-MatchRows 48 49
Swap-MatchRows 59 60
Swap-MatchRows 69 70
Swap-MatchRows 86 87
Swap-MatchRows 118 119
Swap-MatchRows 157 158

# Three-row rotation: 53 <- 54 <- 55 <- 53 (row N takes what used to be
# in row N+1, wrapping around).
$r53 = $ws.Range("F53:V53")
$r54 = $ws.Range("F54:V54")
$r55 = $ws.Range("F55:V55")

$v53 = $r53.Value()
$v54 = $r54.Value()
$v55 = $r55.Value()

$r53.Value = $v54
$r54.Value = $v55
$r55.Value = $v53

# Append the new match as row 182, copying row 181's cell formatting
# (bold/bordered index style in A, date-time number format in E) first.
$ws.Range("A181:V181").Copy()
$ws.Range("A182:V182").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A182").Value = 181
$ws.Range("B182").Value = "italy"
$ws.Range("C182").Value = "serie-a"
$ws.Range("D182").Value = "2023-2024"
$ws.Range("E182").Value = 45296.86458333334
$ws.Range("F182").Value = "Bologna"
$ws.Range("G182").Value = 1
$ws.Range("H182").Value = "Genoa"
$ws.Range("I182").Value = 1
$ws.Range("J182").Value = 1.86
$ws.Range("K182").Value = "23/12/2024 23:02"
$ws.Range("L182").Value = 2.32
$ws.Range("M182").Value = "05/01/2024 20:44"
$ws.Range("N182").Value = 3.34
$ws.Range("O182").Value = "23/12/2024 23:02"
$ws.Range("P182").Value = 2.99
$ws.Range("Q182").Value = "05/01/2024 20:44"
$ws.Range("R182").Value = 4.37
$ws.Range("S182").Value = "23/12/2024 23:02"
$ws.Range("T182").Value = 3.81
$ws.Range("U182").Value = "05/01/2024 20:44"
$ws.Range("V182").Value = "https://www.betexplorer.com/football/italy/serie-a/bologna-genoa/0YBKPhNe/"
